$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: update text to the new, shorter instruction message
$ws.Range("A3").Value = "Email is the only required field."

# Row 4: becomes the header row (First Name / Last Name / Email*)
$ws.Range("A4").Value = "First Name"
$ws.Range("B4").Value = "Last Name"
$ws.Range("C4").Value = "Email*"

# Remove old rows 5-8 (old instructions row 5 & 6, blank spacer row 7, old header row 8)
$ws.Rows("5:8").Delete()

# Row 4 picks up the bold / shaded "header" styling that used to belong to the header row
$headerRow = $ws.Range("A4:Z4")
$headerRow.Font.Bold = $true
$headerRow.Interior.Color = 13421772
